# Apply the "add svgs for params" edit:
#  - Fix the LaTeX strings in the params_deep / params_init sheets by
#    removing the doubled leading backslash (\\alpha -> \alpha, etc.)
#  - Update the active sheet / selection state so that params_init
#    becomes the active tab, with new selections on a few sheets.

$wb = $excel.ActiveWorkbook

$paramsDeep = $wb.Worksheets.Item("params_deep")
$paramsInit = $wb.Worksheets.Item("params_init")
$equations  = $wb.Worksheets.Item("equations")

# --- Fix doubled backslashes in the "latex" column (column B) ---

$paramsDeep.Range("B2").Value = "\alpha"
$paramsDeep.Range("B3").Value = "\beta"
$paramsDeep.Range("B4").Value = "\delta"
$paramsDeep.Range("B5").Value = "\rho"
$paramsDeep.Range("B6").Value = "\rho_{g}"
$paramsDeep.Range("B7").Value = "\rho_{z}"
$paramsDeep.Range("B8").Value = "\xi"
$paramsDeep.Range("B9").Value = "\text{itermax}"

$paramsInit.Range("B2").Value = "\overline{g}"
$paramsInit.Range("B3").Value = "\overline{I_z}"

# --- Update sheet selections / scroll positions ---

[void]$paramsDeep.Range("A7").Select()
[void]$paramsDeep.Range("A9").Select()

[void]$paramsInit.Range("A1").Select()
[void]$paramsInit.Range("B2").Select()

[void]$equations.Range("C6").Select()

# --- Make params_init the active sheet/tab ---

$paramsInit.Activate()
